# Refresh the crypto price/volume snapshot (and restore the original
# relative order of the two Cronos/Binance-PegBSC-USD and
# Kaspa/WhiteBITCoin rows) as produced by the scheduled refresh.
#
# Every assigned value is single-quoted and prefixed with a literal
# leading apostrophe ('). Excel strips that apostrophe and always
# stores the remaining text verbatim as a text/string cell, which is
# what the source workbook uses throughout columns B-E (inline
# strings) -- without it, numeric-looking text like "1.00" or
# "0.810" would silently become the numbers 1 and 0.81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''88.978.28'
$ws.Range("E2").Value = '''  -1.91%  '

$ws.Range("D3").Value = '''3.090.74'
$ws.Range("E3").Value = '''  -3.46%  '

$ws.Range("E4").Value = '''  +0.08%  '

$ws.Range("D5").Value = '''212.77'
$ws.Range("E5").Value = '''  -4.19%  '

$ws.Range("D6").Value = '''623.58'
$ws.Range("E6").Value = '''  -2.79%  '

$ws.Range("D7").Value = '''0.375'
$ws.Range("E7").Value = '''  -6.43%  '

$ws.Range("D8").Value = '''0.810'
$ws.Range("E8").Value = '''  +14.22%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '''  +0.07%  '

$ws.Range("D10").Value = '''3.088.89'
$ws.Range("E10").Value = '''  -3.36%  '

$ws.Range("D11").Value = '''0.614'
$ws.Range("E11").Value = '''  +6.51%  '

$ws.Range("E12").Value = '''  -0.31%  '

$ws.Range("D13").Value = '''0.0000241'
$ws.Range("E13").Value = '''  -7.07%  '

$ws.Range("D14").Value = '''5.29'
$ws.Range("E14").Value = '''  -2.84%  '

$ws.Range("D15").Value = '''88.903.74'
$ws.Range("E15").Value = '''  -1.54%  '

$ws.Range("D16").Value = '''32.27'
$ws.Range("E16").Value = '''  -3.83%  '

$ws.Range("D17").Value = '''3.671.15'

$ws.Range("D18").Value = '''3.095.44'
$ws.Range("E18").Value = '''  -3.29%  '

$ws.Range("E19").Value = '''  +0.75%  '

$ws.Range("D20").Value = '''0.0000211'
$ws.Range("E20").Value = '''  -7.61%  '

$ws.Range("D21").Value = '''13.43'
$ws.Range("E21").Value = '''  -0.35%  '

$ws.Range("D22").Value = '''423.20'
$ws.Range("E22").Value = '''  -3.65%  '

$ws.Range("D23").Value = '''8.27'
$ws.Range("E23").Value = '''  -4.35%  '

$ws.Range("D24").Value = '''4.92'
$ws.Range("E24").Value = '''  -2.89%  '

$ws.Range("D25").Value = '''5.62'
$ws.Range("E25").Value = '''  +4.72%  '

$ws.Range("D26").Value = '''11.90'
$ws.Range("E26").Value = '''  -0.07%  '

$ws.Range("D27").Value = '''82.25'
$ws.Range("E27").Value = '''  +0.99%  '

$ws.Range("D28").Value = '''3.235.31'
$ws.Range("E28").Value = '''  -4.05%  '

$ws.Range("D29").Value = '''1.01'
$ws.Range("E29").Value = '''  +0.74%  '

$ws.Range("B30").Value = '''Cronos'
$ws.Range("C30").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").Value = '''0.172'
$ws.Range("E30").Value = '''  +8.49%  '

$ws.Range("B31").Value = '''Binance-PegBSC-USD'
$ws.Range("C31").Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '''1.08'
$ws.Range("E31").Value = '''  +9.49%  '

$ws.Range("D32").Value = '''8.11'
$ws.Range("E32").Value = '''  -4.28%  '

$ws.Range("D33").Value = '''509.86'
$ws.Range("E33").Value = '''  -5.75%  '

$ws.Range("E34").Value = '''  -12.73%  '

$ws.Range("D35").Value = '''6.75'
$ws.Range("E35").Value = '''  -4.64%  '

$ws.Range("D36").Value = '''1.25'
$ws.Range("E36").Value = '''  -3.87%  '

$ws.Range("D37").Value = '''1.80'
$ws.Range("E37").Value = '''  -6.02%  '

$ws.Range("D38").Value = '''22.28'
$ws.Range("E38").Value = '''  -1.13%  '

$ws.Range("B39").Value = '''Kaspa'
$ws.Range("C39").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.130'
$ws.Range("E39").Value = '''  +2.71%  '

$ws.Range("B40").Value = '''WhiteBITCoin'
$ws.Range("C40").Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '''22.27'
$ws.Range("E40").Value = '''  -0.39%  '

$ws.Range("E41").Value = '''  +0.10%  '

$ws.Range("E42").Value = '''  -0.01%  '

$ws.Range("E43").Value = '''  -2.87%  '

$ws.Range("E44").Value = '''  -6.30%  '

$ws.Range("D45").Value = '''145.98'
$ws.Range("E45").Value = '''  -0.16%  '

$ws.Range("E46").Value = '''  +4.14%  '

$ws.Range("D47").Value = '''0.0693'
$ws.Range("E47").Value = '''  +12.54%  '

$ws.Range("D48").Value = '''43.24'
$ws.Range("E48").Value = '''  -3.69%  '

$ws.Range("D49").Value = '''161.66'
$ws.Range("E49").Value = '''  -6.85%  '

$ws.Range("D50").Value = '''1.21'
$ws.Range("E50").Value = '''  -1.72%  '

$ws.Range("E51").Value = '''  -5.66%  '
